$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append (serial date value, totalScore, C..M = 0, Method = "Random")
$rows = @(
    @{ Row = 9;  Date = 42613.761018518519; Score = 81 },
    @{ Row = 10; Date = 42613.890706018516; Score = 81 },
    @{ Row = 11; Date = 42614.88753472222;  Score = 4  },
    @{ Row = 12; Date = 42615.886805555558; Score = 81 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A: date serial value; copy style from an existing date cell (A8) so the
    # same shared style (numFmtId 22) is reused instead of creating a new style entry.
    $ws.Cells.Item(8, 1).Copy() | Out-Null
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Cells.Item($rowNum, 1).Value = $r.Date

    # Column B: totalScore
    $ws.Cells.Item($rowNum, 2).Value = $r.Score

    # Columns C through M: zero values
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = 0
    }

    # Column N: Method = "Random"
    $ws.Cells.Item($rowNum, 14).Value = "Random"
}
